# Auto-generated edit script applying the Seraph_Profits.xlsx diff
# Updates plain numeric leve-profit-tracking cells across all 8 worksheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR). No formulas are involved -
# every touched cell holds a static number.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1293.8
$ws.Range("I9").Value = 1293.8
$ws.Range("K9").Value = 1293.8
$ws.Range("M9").Value = -1124.8
$ws.Range("H33").Value = 655.5
$ws.Range("I33").Value = 98.05556
$ws.Range("K33").Value = 98.05556
$ws.Range("M33").Value = 130.94444
$ws.Range("H43").Value = 4225
$ws.Range("I43").Value = 1966.6666
$ws.Range("J43").Value = 5580
$ws.Range("K43").Value = 1966.6666
$ws.Range("L43").Value = 5580
$ws.Range("M43").Value = -1897.6666
$ws.Range("N43").Value = -5718
$ws.Range("H62").Value = 6257.364
$ws.Range("I62").Value = 4805.1665
$ws.Range("J62").Value = 8000
$ws.Range("K62").Value = 4805.1665
$ws.Range("L62").Value = 8000
$ws.Range("M62").Value = -4181.1665
$ws.Range("N62").Value = -9248
$ws.Range("H65").Value = 6257.364
$ws.Range("I65").Value = 4805.1665
$ws.Range("J65").Value = 8000
$ws.Range("K65").Value = 24025.8325
$ws.Range("L65").Value = 40000
$ws.Range("M65").Value = -20905.8325
$ws.Range("N65").Value = -46240
$ws.Range("H86").Value = 6000
$ws.Range("J86").Value = 6000
$ws.Range("L86").Value = 6000
$ws.Range("N86").Value = -8246
$ws.Range("H89").Value = 6000
$ws.Range("J89").Value = 6000
$ws.Range("L89").Value = 30000
$ws.Range("N89").Value = -41232
$ws.Range("H103").Value = 3500
$ws.Range("J103").Value = 3500
$ws.Range("L103").Value = 10500
$ws.Range("N103").Value = -11672
$ws.Range("H111").Value = 4246.25
$ws.Range("I111").Value = 3995
$ws.Range("K111").Value = 11985
$ws.Range("M111").Value = -8918
$ws.Range("H112").Value = 2685.8333
$ws.Range("I112").Value = 1605
$ws.Range("J112").Value = 2902
$ws.Range("K112").Value = 4815
$ws.Range("L112").Value = 8706
$ws.Range("M112").Value = -3707
$ws.Range("N112").Value = -10922
$ws.Range("H132").Value = 1418.6842
$ws.Range("I132").Value = 1269.375
$ws.Range("K132").Value = 3808.125
$ws.Range("M132").Value = -1278.125
$ws.Range("H137").Value = 2955.2
$ws.Range("I137").Value = 1741.2142
$ws.Range("J137").Value = 4017.4375
$ws.Range("K137").Value = 5223.642599999999
$ws.Range("L137").Value = 12052.3125
$ws.Range("M137").Value = -2673.642599999999
$ws.Range("N137").Value = -17152.3125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H41").Value = 6854
$ws.Range("J41").Value = 9781
$ws.Range("L41").Value = 9781
$ws.Range("N41").Value = -10609
$ws.Range("H61").Value = 1358
$ws.Range("I61").Value = 1371.6538
$ws.Range("K61").Value = 1371.6538
$ws.Range("M61").Value = -1159.6538
$ws.Range("H74").Value = 3910.7058
$ws.Range("I74").Value = 1978
$ws.Range("K74").Value = 1978
$ws.Range("M74").Value = -1104
$ws.Range("H77").Value = 3910.7058
$ws.Range("I77").Value = 1978
$ws.Range("K77").Value = 9890
$ws.Range("M77").Value = -5522
$ws.Range("H132").Value = 2047.2
$ws.Range("I132").Value = 1681.5555
$ws.Range("K132").Value = 5044.666499999999
$ws.Range("M132").Value = -2514.666499999999
$ws.Range("H136").Value = 1358
$ws.Range("I136").Value = 1371.6538
$ws.Range("K136").Value = 4114.9614
$ws.Range("M136").Value = -1564.9614
$ws.Range("H137").Value = 19999
$ws.Range("I137").Value = 19999
$ws.Range("K137").Value = 19999
$ws.Range("M137").Value = -14899

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H43").Value = 149666.67
$ws.Range("J43").Value = 149666.67
$ws.Range("L43").Value = 149666.67
$ws.Range("N43").Value = -150028.67
$ws.Range("H60").Value = 58710
$ws.Range("J60").Value = 58710
$ws.Range("L60").Value = 58710
$ws.Range("N60").Value = -59908
$ws.Range("H94").Value = 620.7857
$ws.Range("I94").Value = 522.38464
$ws.Range("K94").Value = 522.38464
$ws.Range("M94").Value = -71.38463999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 13151.044
$ws.Range("I99").Value = 10736.462
$ws.Range("K99").Value = 10736.462
$ws.Range("M99").Value = -9238.462
$ws.Range("H126").Value = 13151.044
$ws.Range("I126").Value = 10736.462
$ws.Range("K126").Value = 32209.386
$ws.Range("M126").Value = -29739.386

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 3750
$ws.Range("I70").Value = 2500
$ws.Range("K70").Value = 7500
$ws.Range("M70").Value = -7185
$ws.Range("H73").Value = 3750
$ws.Range("I73").Value = 2500
$ws.Range("K73").Value = 7500
$ws.Range("M73").Value = -6408
$ws.Range("H75").Value = 1107.5
$ws.Range("J75").Value = 1121.75
$ws.Range("L75").Value = 3365.25
$ws.Range("N75").Value = -5361.25
$ws.Range("H78").Value = 1107.5
$ws.Range("J78").Value = 1121.75
$ws.Range("L78").Value = 10095.75
$ws.Range("N78").Value = -20079.75
$ws.Range("H103").Value = 549
$ws.Range("I103").Value = 323.5
$ws.Range("K103").Value = 970.5
$ws.Range("M103").Value = -91.5
$ws.Range("H113").Value = 1914
$ws.Range("J113").Value = 1615.5
$ws.Range("L113").Value = 4846.5
$ws.Range("N113").Value = -9186.5
$ws.Range("H117").Value = 2251.6365
$ws.Range("J117").Value = 3048.1667
$ws.Range("L117").Value = 9144.500100000001
$ws.Range("N117").Value = -16028.5001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3182.0476
$ws.Range("I132").Value = 2447.0557
$ws.Range("K132").Value = 7341.1671
$ws.Range("M132").Value = -4811.1671

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4881.5938
$ws.Range("I132").Value = 3452.9092
$ws.Range("K132").Value = 10358.7276
$ws.Range("M132").Value = -7828.7276

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H132").Value = 1441.92
$ws.Range("I132").Value = 1003
$ws.Range("K132").Value = 3009
$ws.Range("M132").Value = -479
$ws.Range("H136").Value = 44731.668
$ws.Range("I136").Value = 2738.2778
$ws.Range("K136").Value = 8214.8334
$ws.Range("M136").Value = -5664.8334
